$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for the Price column (D) so numeric-looking
# values like "251.06" are not auto-coerced into Number cells -
# they must stay text, matching the rest of the sheet.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "37.402.31"
$ws.Range("E2").Value = "  +4.13%  "

# Row 3
$ws.Range("D3").Value = "2.042.55"
$ws.Range("E3").Value = "  +2.45%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "251.06"
$ws.Range("E5").Value = "  +1.81%  "

# Row 6
$ws.Range("D6").Value = "0.648"
$ws.Range("E6").Value = "  +1.24%  "

# Row 7
$ws.Range("D7").Value = "64.81"
$ws.Range("E7").Value = "  +8.20%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").Value = "0.402"
$ws.Range("E9").Value = "  +9.93%  "

# Row 10
$ws.Range("D10").Value = "59.32"
$ws.Range("E10").Value = "  +1.15%  "

# Row 11
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  +6.62%  "

# Row 12
$ws.Range("E12").Value = "  -0.52%  "

# Row 13
$ws.Range("D13").Value = "0.907"
$ws.Range("E13").Value = "  -4.12%  "

# Row 14
$ws.Range("D14").Value = "23.28"
$ws.Range("E14").Value = "  +20.22%  "

# Row 15
$ws.Range("D15").Value = "14.74"
$ws.Range("E15").Value = "  -0.28%  "

# Row 16
$ws.Range("D16").Value = "2.345.63"
$ws.Range("E16").Value = "  +2.70%  "

# Row 17
$ws.Range("D17").Value = "5.69"
$ws.Range("E17").Value = "  +6.46%  "

# Row 18
$ws.Range("D18").Value = "2.048.94"
$ws.Range("E18").Value = "  +2.92%  "

# Row 19
$ws.Range("D19").Value = "37.305.53"
$ws.Range("E19").Value = "  +4.01%  "

# Row 20
$ws.Range("D20").Value = "72.98"
$ws.Range("E20").Value = "  +1.53%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0879"
$ws.Range("E21").Value = "  +3.14%  "

# Row 22
$ws.Range("D22").Value = "5.46"
$ws.Range("E22").Value = "  +4.53%  "

# Row 23
$ws.Range("D23").Value = "238.78"
$ws.Range("E23").Value = "  +2.24%  "

# Row 24
$ws.Range("E24").Value = "  +0.06%  "

# Row 25
$ws.Range("D25").Value = "2.58"
$ws.Range("E25").Value = "  -2.01%  "

# Row 26
$ws.Range("D26").Value = "2.36"
$ws.Range("E26").Value = "  +3.11%  "

# Row 27
$ws.Range("D27").Value = "9.95"
$ws.Range("E27").Value = "  +3.46%  "

# Row 28
$ws.Range("D28").Value = "160.99"
$ws.Range("E28").Value = "  -2.52%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "19.94"
$ws.Range("E29").Value = "  +2.75%  "

# Row 30
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "0.131"
$ws.Range("E30").Value = "  +33.62%  "

# Row 31
$ws.Range("D31").Value = "0.122"
$ws.Range("E31").Value = "  +2.26%  "

# Row 32
$ws.Range("D32").Value = "5.12"
$ws.Range("E32").Value = "  +3.83%  "

# Row 33
$ws.Range("E33").Value = "  +4.36%  "

# Row 34
$ws.Range("D34").Value = "0.0626"
$ws.Range("E34").Value = "  +3.54%  "

# Row 35
$ws.Range("D35").Value = "4.65"
$ws.Range("E35").Value = "  +5.02%  "

# Row 36
$ws.Range("B36").Value = "THORChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D36").Value = "6.39"
$ws.Range("E36").Value = "  +10.97%  "

# Row 37
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "2.35"
$ws.Range("E37").Value = "  -5.64%  "

# Row 38
$ws.Range("E38").Value = "  +0.24%  "

# Row 39
$ws.Range("E39").Value = "  +2.45%  "

# Row 40
$ws.Range("D40").Value = "2.94"
$ws.Range("E40").Value = "  +26.06%  "

# Row 41
$ws.Range("E41").Value = "  +4.04%  "

# Row 42
$ws.Range("E42").Value = "  +9.38%  "

# Row 43
$ws.Range("E43").Value = "  +6.16%  "

# Row 44
$ws.Range("E44").Value = "  +4.62%  "

# Row 45
$ws.Range("D45").Value = "17.26"
$ws.Range("E45").Value = "  +4.15%  "

# Row 46
$ws.Range("D46").Value = "0.0218"
$ws.Range("E46").Value = "  +2.09%  "

# Row 47
$ws.Range("D47").Value = "95.05"
$ws.Range("E47").Value = "  +0.96%  "

# Row 48
$ws.Range("D48").Value = "7.81"
$ws.Range("E48").Value = "  -0.16%  "

# Row 49
$ws.Range("D49").Value = "1.392.32"
$ws.Range("E49").Value = "  +1.64%  "

# Row 50
$ws.Range("E50").Value = "  +0.38%  "

# Row 51
$ws.Range("D51").Value = "47.01"
$ws.Range("E51").Value = "  +0.12%  "
